$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder/rewrite set-literal text & single values in column E ---
$ws.Range("E18").Value = "{'empty', 'num'}"
$ws.Range("E19").Value = "empty"

$ws.Range("E27").Value = "{'Tuple[None]', 'any'}"
$ws.Range("E28").Value = "Tuple[None]"

$ws.Range("E31").Value = "{'Tuple[None]', 'any'}"
$ws.Range("E32").Value = "Tuple[None]"

$ws.Range("E39").Value = "{'empty', 'num'}"
$ws.Range("E40").Value = "empty"

# --- Fill in the new "Scalpel Accuracy" summary cells on row 211 ---
$ws.Range("C211").Value = "Scalpel Accuracy:"
$ws.Range("D211").Value = 511.76

# --- Move the "Accuracy over PyType" summary down to a new row 212 ---
$ws.Range("E211").Value = ""
$ws.Range("F211").Value = ""

$ws.Range("E212").Value = "Accuracy over PyType"
$ws.Range("F212").Value = 52.94

# Copy the fill/style used throughout the summary rows onto the new row
$ws.Range("A211:F211").Interior.Color = 16777215
$ws.Range("A212:F212").Interior.Color = 16777215
